# Update existing values for May, June, July 2021 rows (174-176)
# and append a new row for August 2021 (row 177).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174 - 01-05-2021
$ws.Cells.Item(174, 2).Value = 63062
$ws.Cells.Item(174, 4).Value = 54013
$ws.Cells.Item(174, 5).Value = 14568
$ws.Cells.Item(174, 7).Value = 12683

# Row 175 - 01-06-2021
$ws.Cells.Item(175, 2).Value = 59052
$ws.Cells.Item(175, 4).Value = 49956

# Row 176 - 01-07-2021
$ws.Cells.Item(176, 2).Value = 56036
$ws.Cells.Item(176, 4).Value = 47162

# Row 177 - new row 01-08-2021
$ws.Cells.Item(177, 1).NumberFormat = "@"
$ws.Cells.Item(177, 1).Value = "01-08-2021"
$ws.Cells.Item(177, 1).Style = "Normal"
$ws.Cells.Item(177, 2).Value = 55844
$ws.Cells.Item(177, 3).Value = 8976
$ws.Cells.Item(177, 4).Value = 46868
$ws.Cells.Item(177, 5).Value = 10852
$ws.Cells.Item(177, 6).Value = 1443
$ws.Cells.Item(177, 7).Value = 9409
